# Paths to Becoming a Penetration Tester.pptx - edit script
#
# Commit: "Add files via upload - Added Certified Ethical Hacker and CompTia
# Security+ certifications for entry level certifications that will help to
# get you through the HR firewall."
#
# This inserts a new slide (title "Decently Respected but low cost certs")
# right after the existing "Free or Low cost Courses" slide (position 7) and
# before the "Offensive Security" slide, pushing every later slide down by
# one position. No other slide content changes.

$p = $ppt.ActivePresentation

# Insert the new slide at position 8 using the standard "Title and Content"
# layout (layout index 2 == ppLayoutText), matching the layout used by all
# of the other certification slides in this deck.
$newSlide = $p.Slides.Add(8, 2)

# --- Title placeholder -------------------------------------------------
$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Decently Respected but low cost certs"
$titleRange.ParagraphFormat.Alignment = 2   # ppAlignCenter

# --- Body / content placeholder ----------------------------------------
$bodyRange = $newSlide.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "Will get you through the HR firewall`rBut don't require a huge amount of experience`rCertified Ethical Hacker: `$1,200 for exam`rCompTia Security+: `$380 for exam"
